$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 392, shifting existing rows 392:410 down to 393:411
$ws.Rows.Item(392).Insert()

# Populate the newly inserted row with the new weekly price record
$ws.Range("A392").Value = 5
$ws.Range("B392").Value = "Macroferia Regional de Talca"
$ws.Range("C392").Value = "Maule"
$ws.Range("D392").Value = 44939
$ws.Range("E392").Value = 7
$ws.Range("F392").Value = 100112003
$ws.Range("G392").Value = "Ajo"
$ws.Range("H392").Value = "Chino"
$ws.Range("I392").Value = "Primera"
$ws.Range("J392").Value = 300
$ws.Range("K392").Value = 17000
$ws.Range("L392").Value = 17000
$ws.Range("M392").Value = 17000
$ws.Range("N392").Value = "$/malla 10 kilos"
$ws.Range("O392").Value = "China"
$ws.Range("P392").Value = 1700
$ws.Range("Q392").Value = 10
$ws.Range("R392").Value = "Hortaliza"
